# Updated cryptos list - applying price/volume/ranking changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.760.56"
$ws.Range("E2").Value = "'  -0.05%  "
$ws.Range("D3").Value = "'1.636.10"
$ws.Range("E3").Value = "'  -0.45%  "
$ws.Range("E4").Value = "'  -0.11%  "
$ws.Range("D5").Value = "'218.05"
$ws.Range("E5").Value = "'  +0.66%  "
$ws.Range("D6").Value = "'0.497"
$ws.Range("E6").Value = "'  -1.12%  "
$ws.Range("E7").Value = "'  -0.06%  "
$ws.Range("D8").Value = "'0.248"
$ws.Range("E8").Value = "'  -1.23%  "
$ws.Range("E9").Value = "'  -0.94%  "
$ws.Range("E10").Value = "'  -1.18%  "
$ws.Range("E11").Value = "'  -0.25%  "
$ws.Range("D12").Value = "'1.862.37"
$ws.Range("E12").Value = "'  -0.53%  "
$ws.Range("D13").Value = "'1.634.60"
$ws.Range("E13").Value = "'  -1.04%  "
$ws.Range("D14").Value = "'4.09"
$ws.Range("E14").Value = "'  -2.71%  "
$ws.Range("D15").Value = "'0.521"
$ws.Range("E15").Value = "'  -1.67%  "
$ws.Range("D16").Value = "'64.09"
$ws.Range("E16").Value = "'  -2.45%  "
$ws.Range("D17").Value = "'26.742.61"
$ws.Range("E17").Value = "'  -0.11%  "
$ws.Range("D18").Value = "'0.0₃0725"
$ws.Range("E18").Value = "'  -2.09%  "
$ws.Range("D19").Value = "'211.41"
$ws.Range("E19").Value = "'  -1.72%  "
$ws.Range("E20").Value = "'  -0.08%  "
$ws.Range("E21").Value = "'  -0.50%  "
$ws.Range("B22").Value = "'Toncoin"
$ws.Range("C22").Value = "'https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D22").Value = "'2.37"
$ws.Range("E22").Value = "'  -3.84%  "
$ws.Range("B23").Value = "'Chainlink"
$ws.Range("C23").Value = "'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").Value = "'6.20"
$ws.Range("E23").Value = "'  -1.63%  "
$ws.Range("D24").Value = "'9.14"
$ws.Range("E24").Value = "'  -2.64%  "
$ws.Range("D25").Value = "'147.93"
$ws.Range("E25").Value = "'  +1.73%  "
$ws.Range("E26").Value = "'  -0.33%  "
$ws.Range("E27").Value = "'  -2.17%  "
$ws.Range("D28").Value = "'6.95"
$ws.Range("E28").Value = "'  -3.24%  "
$ws.Range("D29").Value = "'15.53"
$ws.Range("E29").Value = "'  -1.48%  "
$ws.Range("E30").Value = "'  -3.81%  "
$ws.Range("E31").Value = "'  +0.92%  "
$ws.Range("D32").Value = "'3.35"
$ws.Range("E32").Value = "'  +0.50%  "
$ws.Range("E33").Value = "'  -1.56%  "
$ws.Range("D34").Value = "'1.262.34"
$ws.Range("E34").Value = "'  -0.86%  "
$ws.Range("E35").Value = "'  -1.13%  "
$ws.Range("D36").Value = "'2.45"
$ws.Range("E36").Value = "'  +0.52%  "
$ws.Range("D37").Value = "'0.0173"
$ws.Range("E37").Value = "'  -2.82%  "
$ws.Range("D38").Value = "'0.523"
$ws.Range("E38").Value = "'  -2.65%  "
$ws.Range("E39").Value = "'  +0.05%  "
$ws.Range("D40").Value = "'0.801"
$ws.Range("E40").Value = "'  -3.72%  "
$ws.Range("D41").Value = "'0.801"
$ws.Range("E41").Value = "'  -1.98%  "
$ws.Range("B42").Value = "'MXToken"
$ws.Range("C42").Value = "'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D42").Value = "'2.16"
$ws.Range("E42").Value = "'  -3.82%  "
$ws.Range("B43").Value = "'FraxShare"
$ws.Range("C43").Value = "'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.26"
$ws.Range("E43").Value = "'  -2.06%  "
$ws.Range("B44").Value = "'RocketPoolETH"
$ws.Range("C44").Value = "'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "'1.772.80"
$ws.Range("E44").Value = "'  -1.22%  "
$ws.Range("D45").Value = "'91.59"
$ws.Range("E45").Value = "'  +0.25%  "
$ws.Range("D46").Value = "'59.74"
$ws.Range("E46").Value = "'  +0.93%  "
$ws.Range("D47").Value = "'1.56"
$ws.Range("E47").Value = "'  -2.66%  "
$ws.Range("E48").Value = "'  -0.05%  "
$ws.Range("B49").Value = "'Algorand"
$ws.Range("C49").Value = "'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.0959"
$ws.Range("E49").Value = "'  -2.09%  "
$ws.Range("B50").Value = "'USDD"
$ws.Range("C50").Value = "'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D50").Value = "'1.01"
$ws.Range("E50").Value = "'  -0.15%  "
$ws.Range("E51").Value = "'  -0.58%  "
